# Applies the "Added controllers that can properly handle dips in ph" edit:
#  - Pump 1's "Step" (D2:D6) goes from 120 -> 240
#  - Pump 1's second controller "Step" (J2) goes from 90 -> 700
#  - Updates the sheet's view/selection state to match (topLeftCell B1,
#    active cell I8)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data changes -----------------------------------------------------
$ws.Range("D2").Value = 240
$ws.Range("D3").Value = 240
$ws.Range("D4").Value = 240
$ws.Range("D5").Value = 240
$ws.Range("D6").Value = 240

$ws.Range("J2").Value = 700

# --- View / selection state --------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I8").Select()
